$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.702.10"
$ws.Range("E2").Value = "  +7.60%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.616.06"
$ws.Range("E3").Value = "  +3.89%  "

# Row 4
$ws.Range("E4").Value = "  +0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "418.41"
$ws.Range("E5").Value = "  +0.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.01"
$ws.Range("E6").Value = "  +1.33%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.647"
$ws.Range("E7").Value = "  +3.28%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.608.36"
$ws.Range("E8").Value = "  +3.92%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.768"
$ws.Range("E10").Value = "  +5.51%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.181"
$ws.Range("E11").Value = "  +16.23%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000349"
$ws.Range("E12").Value = "  +55.40%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.67"
$ws.Range("E13").Value = "  +0.32%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.88"
$ws.Range("E14").Value = "  +0.59%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.194.45"
$ws.Range("E15").Value = "  +3.95%  "

# Row 16
$ws.Range("E16").Value = "  -0.28%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.37"
$ws.Range("E17").Value = "  -0.50%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.590.74"
$ws.Range("E18").Value = "  +3.28%  "

# Row 19
$ws.Range("E19").Value = "  +4.42%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.724.12"
$ws.Range("E20").Value = "  +7.65%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.37"
$ws.Range("E21").Value = "  -2.36%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "463.01"
$ws.Range("E22").Value = "  -2.16%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "88.56"
$ws.Range("E23").Value = "  -2.42%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.12"
$ws.Range("E24").Value = "  -5.54%  "

# Row 25
$ws.Range("E25").Value = "  +1.16%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.44"
$ws.Range("E26").Value = "  +3.69%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.01"
$ws.Range("E27").Value = "  -4.34%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.46"
$ws.Range("E28").Value = "  +5.97%  "

# Row 29
$ws.Range("E29").Value = "  +2.02%  "

# Row 30
$ws.Range("E30").Value = "  +3.63%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.36"
$ws.Range("E31").Value = "  +1.80%  "

# Row 32
$ws.Range("E32").Value = "  +3.73%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.38"
$ws.Range("E33").Value = "  -2.25%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.161"
$ws.Range("E34").Value = "  -3.65%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.80"
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.09%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.69"
$ws.Range("E37").Value = "  -2.18%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0493"
$ws.Range("E38").Value = "  +1.01%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0703"
$ws.Range("E39").Value = "  +20.76%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.145"
$ws.Range("E40").Value = "  +7.42%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.14%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.03"
$ws.Range("E42").Value = "  -0.42%  "

# Row 43
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "147.89"
$ws.Range("E43").Value = "  -1.94%  "

# Row 44
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.70"
$ws.Range("E44").Value = "  -4.24%  "

# Row 45
$ws.Range("E45").Value = "  -1.60%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.31"
$ws.Range("E46").Value = "  -2.94%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.308"
$ws.Range("E47").Value = "  -3.92%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.97"
$ws.Range("E48").Value = "  -3.38%  "

# Row 49
$ws.Range("E49").Value = "  -2.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.71"
$ws.Range("E50").Value = "  +16.04%  "

# Row 51
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "115.15"
$ws.Range("E51").Value = "  +5.81%  "
